# Update TestData.xlsx: swap out a few test-data values on Sheet1.
#  - A2: "NEWTC"  -> "Framework_001"
#  - D2: "IE"     -> "Chrome"
#  - D3: "Mozilla"-> "Chrome"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Framework_001"
$ws.Range("D2").Value = "Chrome"
$ws.Range("D3").Value = "Chrome"
